$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "birthday" header in column G, same header style as the other header cells
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").Value = "birthday"

# Birthday sample value + text formatting down through row 3 (even though row 3 stays blank)
$ws.Range("G2:G3").NumberFormat = "@"
$ws.Range("G2").Value = "2005-09-25"

# Validation notes explaining the password behavior, in red to stand out
$ws.Range("H2").Value = "Jika kolom birthday diisi maka password menjadi 20050925"
$ws.Range("H3").Value = "Jika kolom birthday kosong maka password menjadi random"
$ws.Range("H2:H3").Font.Size = 12
$ws.Range("H2:H3").Font.Color = 255

# Update the sample number_id values in column D (kept as text, same style as before)
$ws.Range("D3").Value = "'22394222"
$ws.Range("D2").Value = "'22394213"

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("E7").Select() | Out-Null
